$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.891.28"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.737.06"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.86"
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.97"
$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("E10").Value = "  +4.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("D13").Value = "3.221.76"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("D15").Value = "63.714.30"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "2.742.77"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("E18").Value = "  +1.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "356.28"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.520"
$ws.Range("E23").Value = "  -2.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.20"
$ws.Range("E24").Value = "  -1.55%  "

$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.39"
$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("E29").Value = "  +2.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.39"
$ws.Range("E30").Value = "  +11.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.19"
$ws.Range("E31").Value = "  +1.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.08"
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("E33").Value = "  +0.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").Value = "  +2.67%  "

$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").Value = "  +1.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.981"
$ws.Range("E38").Value = "  -0.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "348.73"
$ws.Range("E39").Value = "  +5.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.28"
$ws.Range("E40").Value = "  +1.84%  "

$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.66"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.88"
$ws.Range("E43").Value = "  +1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.07"
$ws.Range("E44").Value = "  -1.96%  "

$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0252"
$ws.Range("E47").Value = "  -0.95%  "

$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.39"
$ws.Range("E49").Value = "  -1.99%  "

$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  +0.05%  "
